$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 8.898150666666666
$ws.Range("H2").Value = 26.694452
$ws.Range("I2").Value = 0.3765197173862137
$ws.Range("J2").Value = 0.3765197173862137
$ws.Range("M2").Value = 8.908440666666667
$ws.Range("N2").Value = 26.725322
$ws.Range("O2").Value = 0.06231272032629341
$ws.Range("P2").Value = 0.06231272032629341
$ws.Range("Q2").Value = 79.26864725706045
$ws.Range("R2").Value = 713.417825313544
$ws.Range("S2").Value = 0.02346196784682217
$ws.Range("T2").Value = 0.02346196784682217
$ws.Range("G3").Value = 8.898150666666666
$ws.Range("H3").Value = 26.694452
$ws.Range("I3").Value = 0.3765197173862137
$ws.Range("J3").Value = 0.3765197173862137
$ws.Range("O3").Value = 0.3097346304939027
$ws.Range("P3").Value = 0.3097346304939027
$ws.Range("Q3").Value = 394.0165834415849
$ws.Range("R3").Value = 3546.149250974264
$ws.Range("S3").Value = 0.1166211955382876
$ws.Range("T3").Value = 0.1166211955382876
$ws.Range("G4").Value = 8.898150666666666
$ws.Range("H4").Value = 26.694452
$ws.Range("I4").Value = 0.3765197173862137
$ws.Range("J4").Value = 0.3765197173862137
$ws.Range("M4").Value = 89.774269
$ws.Range("N4").Value = 269.322807
$ws.Range("O4").Value = 0.627952649179804
$ws.Range("P4").Value = 0.627952649179804
$ws.Range("Q4").Value = 798.8249715518626
$ws.Range("R4").Value = 7189.424743966764
$ws.Range("S4").Value = 0.236436554001104
$ws.Range("T4").Value = 0.236436554001104
$ws.Range("I5").Value = 0.1415167724465014
$ws.Range("J5").Value = 0.1415167724465015
$ws.Range("M5").Value = 8.908440666666667
$ws.Range("N5").Value = 26.725322
$ws.Range("O5").Value = 0.06231272032629341
$ws.Range("P5").Value = 0.06231272032629341
$ws.Range("Q5").Value = 29.7935077448089
$ws.Range("R5").Value = 268.14156970328
$ws.Range("S5").Value = 0.008818295062938549
$ws.Range("T5").Value = 0.008818295062938551
$ws.Range("I6").Value = 0.1415167724465014
$ws.Range("J6").Value = 0.1415167724465015
$ws.Range("O6").Value = 0.3097346304939027
$ws.Range("P6").Value = 0.3097346304939027
$ws.Range("S6").Value = 0.04383264522240683
$ws.Range("T6").Value = 0.04383264522240684
$ws.Range("I7").Value = 0.1415167724465014
$ws.Range("J7").Value = 0.1415167724465015
$ws.Range("M7").Value = 89.774269
$ws.Range("N7").Value = 269.322807
$ws.Range("O7").Value = 0.627952649179804
$ws.Range("P7").Value = 0.627952649179804
$ws.Range("Q7").Value = 300.2422622338534
$ws.Range("R7").Value = 2702.18036010468
$ws.Range("S7").Value = 0.08886583216115607
$ws.Range("T7").Value = 0.08886583216115608
$ws.Range("G8").Value = 11.39006466666667
$ws.Range("H8").Value = 34.170194
$ws.Range("I8").Value = 0.4819635101672848
$ws.Range("J8").Value = 0.4819635101672848
$ws.Range("M8").Value = 8.908440666666667
$ws.Range("N8").Value = 26.725322
$ws.Range("O8").Value = 0.06231272032629341
$ws.Range("P8").Value = 0.06231272032629341
$ws.Range("Q8").Value = 101.4677152724965
$ws.Range("R8").Value = 913.2094374524681
$ws.Range("S8").Value = 0.03003245741653269
$ws.Range("T8").Value = 0.03003245741653269
$ws.Range("G9").Value = 11.39006466666667
$ws.Range("H9").Value = 34.170194
$ws.Range("I9").Value = 0.4819635101672848
$ws.Range("J9").Value = 0.4819635101672848
$ws.Range("O9").Value = 0.3097346304939027
$ws.Range("P9").Value = 0.3097346304939027
$ws.Range("Q9").Value = 504.360347813701
$ws.Range("R9").Value = 4539.243130323309
$ws.Range("S9").Value = 0.1492807897332083
$ws.Range("T9").Value = 0.1492807897332083
$ws.Range("G10").Value = 11.39006466666667
$ws.Range("H10").Value = 34.170194
$ws.Range("I10").Value = 0.4819635101672848
$ws.Range("J10").Value = 0.4819635101672848
$ws.Range("M10").Value = 89.774269
$ws.Range("N10").Value = 269.322807
$ws.Range("O10").Value = 0.627952649179804
$ws.Range("P10").Value = 0.627952649179804
$ws.Range("Q10").Value = 1022.534729312729
$ws.Range("R10").Value = 9202.812563814559
$ws.Range("S10").Value = 0.3026502630175439
$ws.Range("T10").Value = 0.3026502630175439

$wb.Save()
